# This script reproduces the author's edit: a new data row was inserted
# right after the header/row-14 data block (i.e. before the existing row 15),
# shifting all subsequent rows (old 15..124) down by one (to 16..125).
# The newly inserted row 15 carries its own data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15; this shifts rows 15-124 down to 16-125
# and copies formatting from the row above (row 14), matching the existing
# per-row style (date style on column D, etc.)
$ws.Rows.Item(15).Insert(1)

# Populate the newly inserted row 15 with its data. Most columns repeat the
# constant values used throughout the table; only D, J, K, L, M, P differ
# row-to-row.
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(15, 3).Value = 'Los Lagos'
$ws.Cells.Item(15, 4).Value = (Get-Date -Year 2023 -Month 5 -Day 5).Date
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112026
$ws.Cells.Item(15, 7).Value = 'Haba'
$ws.Cells.Item(15, 8).Value = 'Sin especificar'
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 60
$ws.Cells.Item(15, 11).Value = 23000
$ws.Cells.Item(15, 12).Value = 23000
$ws.Cells.Item(15, 13).Value = 23000
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 920
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = 'Hortaliza'
